# This commit resolves the "cancel ticket" issue:
#  - removes the obsolete "WALLE" sheet
#  - refreshes the login/testing rows on "Details"
#  - updates showtimes/capacity on "Pushpa" and "Krish"
#  - leaves "Pushpa" as the active sheet/tab

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false | Out-Null

$wsDetails = $wb.Worksheets.Item("Details")
$wsPushpa  = $wb.Worksheets.Item("Pushpa")
$wsKrish   = $wb.Worksheets.Item("Krish")
$wsWalle   = $wb.Worksheets.Item("WALLE")

# --- Remove the WALLE sheet entirely ---
$wsWalle.Delete() | Out-Null

# --- "Details" sheet: update existing rows and append new test rows ---
$wsDetails.Range("A2").Value = "hasher"
$wsDetails.Range("A3").Value = "test2"
$wsDetails.Range("A4").Value = "vasa"

$wsDetails.Range("A5").Value = "vasa"
$wsDetails.Range("B5").NumberFormat = "@"
$wsDetails.Range("B5").Value2 = "123"
$wsDetails.Range("B5").Style = "Normal"

$wsDetails.Range("A6").Value = "vasa"
$wsDetails.Range("B6").NumberFormat = "@"
$wsDetails.Range("B6").Value2 = "123"
$wsDetails.Range("B6").Style = "Normal"

$wsDetails.Range("A7").Value = "jerry"
$wsDetails.Range("B7").NumberFormat = "@"
$wsDetails.Range("B7").Value2 = "123"
$wsDetails.Range("B7").Style = "Normal"

# --- "Pushpa" sheet: updated showtimes and capacities ---
$wsPushpa.Range("H2").Value = "10:00-12:00"
$wsPushpa.Range("M2").Value = 100
$wsPushpa.Range("H3").Value = "12:30-02:30"
$wsPushpa.Range("M3").Value = 100
$wsPushpa.Range("H4").Value = "03:00-05:00"
$wsPushpa.Range("M4").Value = 100

# --- "Krish" sheet: updated showtimes and capacities ---
$wsKrish.Range("H2").Value = "10:00-12:00"
$wsKrish.Range("M2").Value = 100
$wsKrish.Range("H3").Value = "12:30-02:30"
$wsKrish.Range("M3").Value = 100
$wsKrish.Range("H4").Value = "03:00-05:00"
$wsKrish.Range("M4").Value = 100

# --- Restore each sheet's selection, then leave Pushpa as the active tab ---
$wsDetails.Activate() | Out-Null
$wsDetails.Range("D7").Select() | Out-Null

$wsKrish.Activate() | Out-Null
$wsKrish.Range("M11").Select() | Out-Null

$wsPushpa.Activate() | Out-Null
$wsPushpa.Range("M2:M4").Select() | Out-Null
